$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the raw (non-formula) hour entries for "Donderdag" (Thursday) of Week 11 (row 53)
$ws.Range("C53:H53").Value = 2
$ws.Range("I53").Value = 1

# Update the manually entered "Totaal Game-Lab uren p/w" total for Week 11 (row 55)
$ws.Range("B55").Value = 11

# Recalculate all dependent formulas (K2, L2:L8, M2:M10, C55:I55, etc.)
$excel.Calculate()

# Restore the view state: active cell selection
$ws.Activate()
$ws.Range("K49").Select()
